$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 get cyclically shifted down by one: new row N = old row (N-1),
# with old row 5 wrapping around to become the new row 2.
$data = @(
    @(5,7,3,6,2,8,1,4),
    @(2,4,1,6,3,7,5,8),
    @(3,8,2,7,4,6,1,5),
    @(6,7,1,8,4,5,2,3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}

$ws.Range("A2:H5").Select()
